# The edit reorders the data rows (rows 2-33) of the sheet: each
# destination row ends up holding the values that used to live in a
# different source row (row 25 is the only row that stays put). Column
# layout (A:T) is unchanged - only which physical row holds which
# record changes. We snapshot every row's values first (so the
# permutation can't clobber a source row before it has been read),
# then write the snapshots back out in their new positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row number -> source row number (both in the
# original, pre-edit layout).
$map = @{
    2  = 29
    3  = 2
    4  = 20
    5  = 33
    6  = 12
    7  = 11
    8  = 28
    9  = 7
    10 = 26
    11 = 13
    12 = 6
    13 = 4
    14 = 18
    15 = 24
    16 = 17
    17 = 9
    18 = 10
    19 = 5
    20 = 22
    21 = 8
    22 = 15
    23 = 16
    24 = 30
    25 = 25
    26 = 32
    27 = 21
    28 = 31
    29 = 14
    30 = 3
    31 = 27
    32 = 23
    33 = 19
}

# Snapshot all source rows (2-33) before writing anything back.
$snapshot = @{}
foreach ($r in 2..33) {
    $snapshot[$r] = $ws.Range("A$r`:T$r").Value2
}

# Write each destination row from its mapped source snapshot.
foreach ($destRow in 2..33) {
    $srcRow = $map[$destRow]
    $ws.Range("A$destRow`:T$destRow").Value2 = $snapshot[$srcRow]
}
